$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.470.95'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '1.803.87'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.57'
$ws.Range('E5').Value = '  -1.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.602'
$ws.Range('E6').Value = '  +4.19%  '
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '39.66'
$ws.Range('E8').Value = '  +7.66%  '
$ws.Range('E9').Value = '  -3.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0669'
$ws.Range('E10').Value = '  -3.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0983'
$ws.Range('E11').Value = '  +1.92%  '
$ws.Range('D12').Value = '2.066.51'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.99'
$ws.Range('E13').Value = '  -5.26%  '
$ws.Range('D14').Value = '1.800.27'
$ws.Range('E14').Value = '  -0.24%  '
$ws.Range('E15').Value = '  -3.38%  '
$ws.Range('D16').Value = '34.457.49'
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.37'
$ws.Range('E17').Value = '  -2.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.14'
$ws.Range('E18').Value = '  -2.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '239.54'
$ws.Range('E19').Value = '  -2.47%  '
$ws.Range('D20').Value = '0.0₃0768'
$ws.Range('E20').Value = '  -3.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.14'
$ws.Range('E21').Value = '  -4.13%  '
$ws.Range('E22').Value = '  -0.30%  '
$ws.Range('E23').Value = '  -2.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.18'
$ws.Range('E24').Value = '  -0.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.52'
$ws.Range('E25').Value = '  -0.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '17.68'
$ws.Range('E26').Value = '  +4.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.68'
$ws.Range('E27').Value = '  -4.07%  '
$ws.Range('E28').Value = '  +2.92%  '
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('E30').Value = '  -1.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.75'
$ws.Range('E31').Value = '  -2.74%  '
$ws.Range('E32').Value = '  -3.28%  '
$ws.Range('E33').Value = '  -4.53%  '
$ws.Range('E34').Value = '  +0.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.641'
$ws.Range('E35').Value = '  -4.94%  '
$ws.Range('D36').Value = '1.302.71'
$ws.Range('E36').Value = '  -6.75%  '
$ws.Range('E37').Value = '  -1.49%  '
$ws.Range('E38').Value = '  -2.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.31'
$ws.Range('E39').Value = '  -6.71%  '
$ws.Range('E40').Value = '  +0.21%  '
$ws.Range('E41').Value = '  +1.74%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '81.87'
$ws.Range('E42').Value = '  -1.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.949'
$ws.Range('E43').Value = '  -1.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.79'
$ws.Range('E44').Value = '  -1.66%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.92'
$ws.Range('E45').Value = '  +2.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0515'
$ws.Range('E46').Value = '  +3.18%  '
$ws.Range('D47').Value = '1.966.95'
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.73'
$ws.Range('E48').Value = '  -5.23%  '
$ws.Range('E49').Value = '  -0.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '102.13'
$ws.Range('E50').Value = '  -2.51%  '
$ws.Range('E51').Value = '  -0.34%  '
